$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.182.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("D3").Value = '''2.415.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.38%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''561.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.50%  '
$ws.Range("D6").Value = '''142.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.63%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '''0.530'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.12%  '
$ws.Range("D9").Value = '''2.410.92'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("E10").Value = '  +1.99%  '
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").Value = '''0.352'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").Value = '''25.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("D15").Value = '''0.0000175'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("D16").Value = '''2.851.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("D17").Value = '''62.091.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.29%  '
$ws.Range("D18").Value = '''2.411.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.24%  '
$ws.Range("D19").Value = '''11.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.78%  '
$ws.Range("D20").Value = '''323.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("D22").Value = '''6.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.93%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '''65.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("D25").Value = '''1.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").Value = '''9.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.62%  '
$ws.Range("D27").Value = '''576.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.32%  '
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("E29").Value = '  +2.27%  '
$ws.Range("D30").Value = '''0.0₃0944'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.88%  '
$ws.Range("D31").Value = '''8.22'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("E32").Value = '  +3.47%  '
$ws.Range("E33").Value = '  +1.80%  '
$ws.Range("E34").Value = '  +1.72%  '
$ws.Range("E35").Value = '  +1.87%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '''5.54'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.19%  '
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("D39").Value = '''152.89'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.56%  '
$ws.Range("D40").Value = '''0.382'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("E42").Value = '  -4.37%  '
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").Value = '''2.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.03%  '
$ws.Range("D45").Value = '''148.53'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.16%  '
$ws.Range("D46").Value = '''3.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("E47").Value = '  +1.48%  '
$ws.Range("D48").Value = '''20.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("E49").Value = '  +2.50%  '
$ws.Range("D50").Value = '''0.0917'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.84%  '
